# Sample-Project-Report.xlsx edit script
# Adds a new "Project_Summary" column (I) to the ProjectReport sheet,
# extends the header formatting / column width, updates the AutoFilter
# and the _FilterDatabase defined name to cover the new column, and
# updates the saved selection / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new header cell in I1 --------------------------------------
$ws.Range("I1").Value = "Project_Summary"

# Copy the header formatting (font/style) from H1 onto the new I1 cell so
# it keeps the same bold blue "Segoe UI" header style (cellXf s="2").
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Give column I a sensible width similar to the other bestFit columns.
# (The runtime quantises ColumnWidth to 1/6-character increments, so we pick
# the input that rounds to the closest achievable width to 15.81640625.)
$ws.Columns.Item(9).ColumnWidth = 14.95

# --- Re-apply the AutoFilter across the new range A1:I1 ------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:I1").AutoFilter()

# --- Update the hidden _FilterDatabase defined name to match -------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ProjectReport!_FilterDatabase") {
        $n.RefersTo = "=ProjectReport!`$A`$1:`$I`$1"
    }
}

# --- Update the saved selection / scroll position -------------------------
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E12").Select()
